# Updated data from NHC, add single point estimate according to suspected number
#
# This script reproduces (as far as the COM surface allows) the authoring
# session captured by the target diff:
#   1. A Page Layout > Colors > Customize Colors edit that changed the
#      theme's "Background 1" (lt1) swatch from white to a pale green
#      (CCE8CF). The runtime's Theme object model is read-only for
#      serialization purposes, so this call is best-effort / a no-op on
#      save, but it documents the intended edit and is harmless to issue.
#   2. Refreshed Confirmed/Suspected case counts for 2020-01-19 .. 2020-01-27
#      (rows 19-27), including newly-added Suspected figures for rows 24-27.
#   3. The active selection left on B19 (matching the last-edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# --- 1. Theme "Background 1" color customization -----------------------
try {
    $wb.Theme.ThemeColorScheme.Colors(2) = RGB(204, 232, 207)   # CCE8CF
} catch {
    # Theme editing may be unavailable in this environment; continue.
}

# --- 2. Updated Confirmed (B) / Suspected (C) counts --------------------
$ws.Range("B19").Value = 62
$ws.Range("B20").Value = 121
$ws.Range("B21").Value = 198
$ws.Range("B22").Value = 291
$ws.Range("B23").Value = 440
$ws.Range("B24").Value = 571
$ws.Range("C24").Value = 393
$ws.Range("B25").Value = 830
$ws.Range("C25").Value = 1072
$ws.Range("B26").Value = 1287
$ws.Range("C26").Value = 1965
$ws.Range("B27").Value = 1975
$ws.Range("C27").Value = 2684

# --- 3. Leave selection on the last-touched cell -------------------------
$ws.Range("B19").Select() | Out-Null
